$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weights")

$ws.Range("B4").Value = 4.307709026009078
$ws.Range("C4").Value = -12.04447704155604
$ws.Range("D4").Value = 8.26792402307002
$ws.Range("E4").Value = -14.85089138709381

$ws.Range("B5").Value = 3.767339116572964
$ws.Range("C5").Value = -11.667391693184
$ws.Range("D5").Value = 9.064465412192149
$ws.Range("E5").Value = -12.80714756008984

$ws.Range("B6").Value = -2
$ws.Range("C6").Value = -6
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = -5
